# Generate Report for Handoff
# Updates the localization-status workbook: marks the "b.md" files as
# "Ready for handoff" (handoff re-issued because handback was stale) for
# both the zh-cn and de-de language sheets, as well as the Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a4c0eb9560ded1383d98fbcaa0523bbe689564c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9c06fb710b551939813ca1e46cc25360f636e3fa/e2e/b.md."

# --- Overview sheet: row 3 is the b.md entry ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-16 14:37:49"

# --- zh-cn sheet: row 3 is the b.md entry ---
# (ColumnWidth is in "characters"; the engine adds a fixed 5/6 padding when
# serialising to the OOXML <col width=".."> attribute, so back it off here
# to land exactly on width="40" in the saved file.)
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces the literal "False" to be stored as text rather
# than auto-coerced into a Boolean; reapply the Normal style afterwards so
# the quote-prefix flag that the apostrophe trick sets doesn't linger on
# the cell's format.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-16 14:37:44"
$zhcn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 is the b.md entry ---
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-16 14:37:49"
$dede.Range("P3").Value = $errorDetail
